# Blog Statistics - November 2016
# Fill in the November (column M) figure for year 2016 (row 13), which was
# previously a blank placeholder ("-"), with the actual hit count. Downstream
# formulas (row Total P13, Grand-total P14, and the trend chart's cached
# series) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Enter the November 2016 figure.
$ws.Range("M13").Value = 713080

# Match the number formatting/style already used by the other filled-in
# month cells in this row (I13:L13) instead of the "blank" placeholder style.
$ws.Range("L13").Copy()
$ws.Range("M13").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Reflect the latest cell the user worked on/selected.
[void]$ws.Range("S8").Select()
